$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")

# --- Clear the stray, formatted-but-empty C233 / C234 cells ---
$ws.Range("C233:C234").Clear()

# --- Add the new variable rows (235-243) ---
# Use an existing, correctly-formatted row (232) as the formatting template
# for column B (style 4) and column C (style 5) so we reuse the workbook's
# existing style entries instead of creating new ones.

$newRows = @(
    @{ Row = 235; A = 232; B = "Diff_MPuse_g";      C = "calculate_Body_MPuse_g_Trg"; D = "No"; E = 2689; PasteTo = "E" },
    @{ Row = 236; A = 233; B = "An_MEIn_approx";    C = "calculate_MP_requirement";   D = "No"; E = 2686; PasteTo = "E" },
    @{ Row = 237; A = 234; B = "Frm_MPUse_g_Trg";   C = "calculate_MP_requirement";   D = "No"; E = $null; PasteTo = "E" },
    @{ Row = 238; A = 235; B = "Kg_MP_NP_Trg";      C = "calculate_MP_requirement";   D = "No"; E = $null; PasteTo = "E" },
    @{ Row = 239; A = 236; B = "Min_MPuse_g";       C = "calculate_MP_requirement";   D = "No"; E = $null; PasteTo = "E" },
    @{ Row = 240; A = 237; B = "Frm_NPgain_g";      C = "calculate_MP_requirement";   D = "No"; E = $null; PasteTo = "E" },
    @{ Row = 241; A = 238; B = "Kg_MP_NP_Trg";      C = "calculate_MP_requirement";   D = "No"; E = $null; PasteTo = "D" },
    @{ Row = 242; A = 239; B = "Rsrv_NPgain_g";     C = "calculate_MP_requirement";   D = "No"; E = $null; PasteTo = "D" },
    @{ Row = 243; A = 240; B = "Rsrv_MPUse_g_Trg";  C = "calculate_MP_requirement";   D = "No"; E = $null; PasteTo = "D" }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row

    # Copy the number/text formatting from row 232 (through column E or D,
    # depending on how wide the author's original paste was for this row)
    # into the new row first, then overwrite the values - this keeps column
    # B on style 4 and column C on style 5, matching the rest of the table.
    $ws.Range("B232:" + $r.PasteTo + "232").Copy()
    $ws.Range("B" + $rowNum + ":" + $r.PasteTo + $rowNum).PasteSpecial(-4122)

    $ws.Cells.Item($rowNum, 1).Value = $r.A
    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D
    if ($r.E -ne $null) {
        $ws.Cells.Item($rowNum, 5).Value = $r.E
    } elseif ($r.PasteTo -eq "E") {
        # Clear any pasted formatting in column E that isn't backed by a
        # real value, while keeping the row's recorded column extent (the
        # author's rows 237-240 still carry a 1:6 span even with no E/F
        # cell present).
        $ws.Range("E" + $rowNum).Clear()
    }
}

$excel.CutCopyMode = 0

# --- Update the view state to match where the author ended up ---
$ws.Activate()
$ws.Range("D237:D243").Select()
